$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: L3's year value moves from 2020 to 2021, and a new column M
# is appended with the next year (2022), matching L3's formatting.
$ws.Range("L3").Copy()
$ws.Range("M3").PasteSpecial(-4122)
$ws.Range("M3").Value = 2022
$ws.Range("L3").Value = 2021

# Row 4: extend the data series into the new column M, matching L4's
# formatting (same indicator value carried over, 6.18).
$ws.Range("L4").Copy()
$ws.Range("M4").PasteSpecial(-4122)
$ws.Range("M4").Value = 6.18

$excel.CutCopyMode = $false

# The active selection moves to M9 after the edit.
$ws.Range("M9").Select()
